$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the fractional precision of the existing I39 timestamp value
$ws.Range("I39").Value = 45521.98183673611

# Capture the date/time number format used by the "order_date" column so
# the newly appended rows keep the same display format (style index 2).
$dateFormat = $ws.Range("I39").NumberFormat

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "Коммутатор"
$ws.Range("C40").Value = "ПРГС.465000.028"
$ws.Range("D40").Value = 975
$ws.Range("E40").Value = 984
$ws.Range("F40").Value = "10.8.15.6"
$ws.Range("G40").Value = "10.8.15.15"
$ws.Range("H40").Value = 10
$ws.Range("I40").Value = 45581.50276196759
$ws.Range("I40").NumberFormat = $dateFormat

# Row 41
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "Монитор"
$ws.Range("C41").Value = "ПРГС.465000.012"
$ws.Range("D41").Value = 1025
$ws.Range("E41").Value = 1034
$ws.Range("F41").Value = "10.8.99.0"
$ws.Range("G41").Value = "10.8.99.9"
$ws.Range("H41").Value = 10
$ws.Range("I41").Value = 45581.54614268518
$ws.Range("I41").NumberFormat = $dateFormat

# Row 42
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "Монитор"
$ws.Range("C42").Value = "ПРГС.465000.012"
$ws.Range("D42").Value = 1035
$ws.Range("E42").Value = 1044
$ws.Range("F42").Value = "10.8.99.10"
$ws.Range("G42").Value = "10.8.99.19"
$ws.Range("H42").Value = 10
$ws.Range("I42").Value = 45581.54656073602
$ws.Range("I42").NumberFormat = $dateFormat
